# Gaussian Quadrature Scheme export: append the next averaged-intensity
# row (HKL index 14, "HexGrid-60degTilt5degRes") to the bottom of the
# table on Sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 15 (A15:M15) carries the formatting used for every row label in
# column A (bold font, thin border, centered/top aligned). Copy that
# formatting down to A16 before writing the new values so the new cell
# reuses the existing style instead of minting a near-duplicate one.
$ws.Range("A15").Copy()
$ws.Range("A16").PasteSpecial(-4122)

$ws.Range("A16").Value = 14

# Column B reuses the same scheme label as row 15 ("HexGrid-60degTilt5degRes"),
# so copy its text across rather than retyping it (keeps the same shared
# string, bumping sharedStrings' count but not uniqueCount).
$ws.Range("B16").Value = $ws.Range("B15").Text

$ws.Range("C16").Value = 0.9671164926603516
$ws.Range("D16").Value = 1.1173012170782
$ws.Range("E16").Value = 0.9646396893408864
$ws.Range("F16").Value = 0.9671164926603516
$ws.Range("G16").Value = 1.06341495830896
$ws.Range("H16").Value = 0.9131063717481164
$ws.Range("I16").Value = 0.9658307839933642
$ws.Range("J16").Value = 1.1173012170782
$ws.Range("K16").Value = 1.040970453209543
$ws.Range("L16").Value = 1.004043472934947
$ws.Range("M16").Value = 0.9985682521883131
